$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("example-slideshow")

# Add a video clip (jump-07-0.mp4) right before the jump-07.jpg row (row 9
# in the not-yet-shifted sheet). This will end up at final row 10 once the
# audio row below is inserted at the top.
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 2).Value = "jump-07-0.mp4"

# Add an audio clip (ding.mp3) as the new first data row, with an option
# that ends the trial once the audio finishes playing.
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 2).Value = "ding.mp3"

# Give the jump-03.jpg row (now final row 5) an explicit "image" type.
$ws.Cells.Item(5, 1).Value = "image"

# Set the options JSON for the new audio row.
$ws.Cells.Item(2, 3).Value = '{"trial_ends_after_audio":true}'

# Match the authored selection state.
$ws.Range("C3").Select()
